# October 11th Work. (#83)
# Populates the daily-progress log rows 12-21 (Oct 11 - Oct 20, 2025) with the
# skills/backgrounds/talents rolled that day, reusing the same three
# conditional-format cell styles (Good/Neutral/Bad) already used on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Write the cell text/values -----------------------------------------

# Row 12
$ws.Range('B12').Value = 'Sailor'
$ws.Range('C12').Value = 'Vitalité'
$ws.Range('E12').Value = 'Dur à cuire, Effacé'
$ws.Range('F12').Value = 'Sous contrainte'
$ws.Range('G12').Value = 'Diplomatie'
$ws.Range('H12').Value = 'Herboristerie'
$ws.Range('I12').Value = 'Lutte'
$ws.Range('J12').Value = 'Cambrioleur'

# Row 13
$ws.Range('B13').Value = 'Soldier'
$ws.Range('C13').Value = 'Artisanat'
$ws.Range('E13').Value = 'Faiblesse, Féroce'
$ws.Range('F13').Value = 'Artisan'
$ws.Range('G13').Value = 'Discipline'
$ws.Range('H13').Value = 'Initiative accrue'
$ws.Range('J13').Value = 'Catalyste'

# Row 14
$ws.Range('B14').Value = 'Urchin'
$ws.Range('C14').Value = 'Athlétisme'
$ws.Range('E14').Value = 'Entrepôt à connaissances'
$ws.Range('F14').Value = 'Bohème'
$ws.Range('G14').Value = 'Dressage'
$ws.Range('J14').Value = 'Chaman'

# Row 15
$ws.Range('B15').Value = 'Other Backgrounds'
$ws.Range('C15').Value = 'Connaissance'
$ws.Range('E15').Value = 'Fignolage, Globe-trotteur'
$ws.Range('F15').Value = 'Exilé'
$ws.Range('G15').Value = 'Élémentarisme'

# Row 16
$ws.Range('B16').Value = 'Other Backgrounds'
$ws.Range('C16').Value = 'Diplomatie'
$ws.Range('F16').Value = 'Guérisseur'

# Row 17
$ws.Range('B17').Value = 'Other Backgrounds'
$ws.Range('C17').Value = 'Discipline'
$ws.Range('F17').Value = 'Marchand'

# Row 18
$ws.Range('B18').Value = 'Other Backgrounds'
$ws.Range('C18').Value = 'Furtivité'
$ws.Range('F18').Value = 'Milicien'

# Row 19
$ws.Range('C19').Value = 'Intuition'
$ws.Range('F19').Value = 'Noble'

# Row 20
$ws.Range('C20').Value = 'Investigation'
$ws.Range('F20').Value = 'Paysan'

# Row 21
$ws.Range('C21').Value = 'Linguistique'
$ws.Range('F21').Value = 'Religieux'

# --- 2) Re-apply the existing colored cell styles ---------------------------
# Row 2 already carries the three fill styles used across the sheet:
#   B2 -> "Neutral" (yellow, style index 6)
#   C2 -> "Good"    (green,  style index 5)
#   D2 -> "Bad"     (red,    style index 7)
# Copy format-only (xlPasteFormats = -4122) from those donor cells onto the
# new cells below, so the same shared cell-style indices get reused instead of
# Excel minting new ones.

$cellsStyle5 = @('C12', 'E12', 'G12', 'H12', 'I12', 'G13', 'H13', 'G14', 'G15')
$ws.Range('C2').Copy() | Out-Null  # style 5 = "Good"
foreach ($ref in $cellsStyle5) {
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

$cellsStyle6 = @('B12', 'F12', 'J12', 'B13', 'C13', 'E13', 'I13', 'J13', 'B14', 'C14', 'E14', 'H14', 'I14', 'J14', 'B15', 'C15', 'E15', 'B16', 'C16', 'B17', 'C17', 'B18', 'C18', 'B19', 'C19', 'B20', 'C20', 'B21', 'C21')
$ws.Range('B2').Copy() | Out-Null  # style 6 = "Neutral"
foreach ($ref in $cellsStyle6) {
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

$cellsStyle7 = @('D12', 'K12', 'L12', 'M12', 'D13', 'F13', 'K13', 'L13', 'M13', 'D14', 'F14', 'K14', 'L14', 'M14', 'D15', 'F15', 'K15', 'L15', 'M15', 'D16', 'E16', 'F16', 'K16', 'L16', 'M16', 'D17', 'E17', 'F17', 'K17', 'L17', 'M17', 'D18', 'E18', 'F18', 'K18', 'L18', 'M18', 'D19', 'E19', 'F19', 'K19', 'L19', 'M19', 'D20', 'E20', 'F20', 'K20', 'L20', 'M20', 'D21', 'E21', 'F21', 'K21', 'L21', 'M21')
$ws.Range('D2').Copy() | Out-Null  # style 7 = "Bad"
foreach ($ref in $cellsStyle7) {
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

Write-Host "Rows 12-21 populated."